$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Duplicate the "osu! writeup" sheet so the original keeps its sheetId/rId
# (it will remain in slot 2, renamed to "Sheet3"), while the new copy
# (fresh sheetId, placed right after) becomes the "real" osu! writeup sheet.
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$newOsu = $wb.Worksheets.Item(3)
$newOsu.Name = "osu! writeup temp"

# Clear out the old osu! writeup content/layout from the original sheet
# object (now destined to become "Sheet3"), then repopulate it with a copy
# of the first 14 rows of "ep 1".
$ws2.Cells.Clear()

$ws1.Range("A1:E1").Copy($ws2.Range("A1:E1"))
$ws1.Range("A2:C2").Copy($ws2.Range("A2:C2"))
$ws1.Range("A3:C3").Copy($ws2.Range("A3:C3"))
$ws1.Range("A4:C4").Copy($ws2.Range("A4:C4"))
$ws1.Range("A5:C5").Copy($ws2.Range("A5:C5"))
$ws1.Range("A6:C6").Copy($ws2.Range("A6:C6"))
$ws1.Range("A7:C7").Copy($ws2.Range("A7:C7"))
$ws1.Range("A8:C8").Copy($ws2.Range("A8:C8"))
$ws1.Range("A9:E9").Copy($ws2.Range("A9:E9"))
$ws1.Range("A10:C10").Copy($ws2.Range("A10:C10"))
$ws1.Range("A11:C11").Copy($ws2.Range("A11:C11"))
$ws1.Range("A12:C12").Copy($ws2.Range("A12:C12"))
$ws1.Range("A13:C13").Copy($ws2.Range("A13:C13"))
$ws1.Range("A14:C14").Copy($ws2.Range("A14:C14"))

$ws2.Name = "Sheet3"
$newOsu.Name = "osu! writeup"
